$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 already has the underline style (s="1") from the original file but was empty.
# Fill it in with "EUR", keeping its existing style.
$ws.Range("A4").Value = "EUR"

# Add two new rows with default style.
$ws.Range("A5").Value = "JPY"
$ws.Range("A6").Value = "AUD"

# Move the active selection to A6, matching the post-edit state.
$ws.Range("A6").Select()
